$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name labels (column A) for rows whose displayed country changed ---
$ws.Range("A1").Value = 'Datos actualizados a 31 de Marzo de 2020 a las 13:50'
$ws.Range("A143").Value = 'Mali'
$ws.Range("A144").Value = 'Niger'
$ws.Range("A145").Value = 'Bermudas'
$ws.Range("A146").Value = 'Etiopia'
$ws.Range("A167").Value = 'Siria'
$ws.Range("A168").Value = 'Groenlandia'
$ws.Range("A169").Value = 'Suazilandia'
$ws.Range("A170").Value = 'Laos'
$ws.Range("A171").Value = 'Granada'
$ws.Range("A174").Value = 'San Cristobal y Nieves'
$ws.Range("A177").Value = 'Mozambique'
$ws.Range("A179").Value = 'Guyana'
$ws.Range("A180").Value = 'Antigua y Barbuda'

# --- Update numeric statistic cells (columns B-H) ---
$ws.Range("D8").Value = 15824
$ws.Range("E8").Value = 50576
$ws.Range("G8").Value = 6
$ws.Range("H8").Value = 651
$ws.Range("B19").Value = 7443
$ws.Range("C19").Value = 1035
$ws.Range("E19").Value = 7240
$ws.Range("F19").Value = 188
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 160
$ws.Range("D92").Value = 57
$ws.Range("E92").Value = 150
$ws.Range("B108").Value = 132
$ws.Range("C108").Value = 10
$ws.Range("E108").Value = 114
$ws.Range("B143").Value = 28
$ws.Range("C143").Value = 3
$ws.Range("E143").Value = 26
$ws.Range("H143").Value = 2
$ws.Range("D144").Value = 0
$ws.Range("E144").Value = 24
$ws.Range("H144").Value = 3
$ws.Range("B145").Value = 27
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 10
$ws.Range("E145").Value = 17
$ws.Range("F145").Value = 0
$ws.Range("C146").Value = 2
$ws.Range("D146").Value = 2
$ws.Range("F146").Value = 2
$ws.Range("H146").Value = 0
$ws.Range("D167").Value = 0
$ws.Range("H167").Value = 2
$ws.Range("D168").Value = 2
$ws.Range("H168").Value = 0
$ws.Range("C170").Value = 1
$ws.Range("C171").Value = 0
$ws.Range("C174").Value = 1
$ws.Range("E177").Value = 8
$ws.Range("H177").Value = 0
$ws.Range("B179").Value = 8
$ws.Range("H179").Value = 1
